$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $oldStyle = $range.Style
    $range.Value = "'" + $value
    $range.Style = $oldStyle
}

Set-TextValue $ws.Range("D2") "30.200.99"
Set-TextValue $ws.Range("E2") "  +3.24%  "

Set-TextValue $ws.Range("D3") "1.896.54"
Set-TextValue $ws.Range("E3") "  +0.17%  "

Set-TextValue $ws.Range("E4") "  -0.34%  "

Set-TextValue $ws.Range("D5") "325.23"

Set-TextValue $ws.Range("E6") "  -0.33%  "

Set-TextValue $ws.Range("D7") "0.5165"
Set-TextValue $ws.Range("E7") "  +0.57%  "

Set-TextValue $ws.Range("D8") "0.4013"
Set-TextValue $ws.Range("E8") "  +2.74%  "

Set-TextValue $ws.Range("D9") "0.08441"
Set-TextValue $ws.Range("E9") "  +0.23%  "

Set-TextValue $ws.Range("D10") "42.69"
Set-TextValue $ws.Range("E10") "  +0.67%  "

Set-TextValue $ws.Range("E11") "  +0.48%  "

Set-TextValue $ws.Range("D12") "23.16"
Set-TextValue $ws.Range("E12") "  +12.28%  "

Set-TextValue $ws.Range("D13") "6.430"
Set-TextValue $ws.Range("E13") "  +3.17%  "

Set-TextValue $ws.Range("D14") "1.899.84"
Set-TextValue $ws.Range("E14") "  +0.62%  "

Set-TextValue $ws.Range("D15") "7.337"
Set-TextValue $ws.Range("E15") "  +0.37%  "

Set-TextValue $ws.Range("E16") "  -0.33%  "

Set-TextValue $ws.Range("D17") "94.75"
Set-TextValue $ws.Range("E17") "  +2.00%  "

Set-TextValue $ws.Range("D18") "0.00001109"
Set-TextValue $ws.Range("E18") "  +0.45%  "

Set-TextValue $ws.Range("D19") "0.06667"
Set-TextValue $ws.Range("E19") "  -1.10%  "

Set-TextValue $ws.Range("D20") "18.26"
Set-TextValue $ws.Range("E20") "  +2.54%  "

Set-TextValue $ws.Range("E21") "  -0.26%  "

Set-TextValue $ws.Range("D22") "5.944"
Set-TextValue $ws.Range("E22") "  -0.97%  "

Set-TextValue $ws.Range("D23") "30.211.06"
Set-TextValue $ws.Range("E23") "  +3.24%  "

Set-TextValue $ws.Range("D24") "11.30"
Set-TextValue $ws.Range("E24") "  +1.68%  "

Set-TextValue $ws.Range("D25") "2.211"
Set-TextValue $ws.Range("E25") "  -0.07%  "

Set-TextValue $ws.Range("D26") "2.114.25"
Set-TextValue $ws.Range("E26") "  +0.34%  "

Set-TextValue $ws.Range("D27") "21.72"
Set-TextValue $ws.Range("E27") "  +4.19%  "

Set-TextValue $ws.Range("D28") "161.20"
Set-TextValue $ws.Range("E28") "  +1.21%  "

Set-TextValue $ws.Range("D29") "2.381"
Set-TextValue $ws.Range("E29") "  -1.82%  "

Set-TextValue $ws.Range("D30") "129.03"
Set-TextValue $ws.Range("E30") "  +1.26%  "

Set-TextValue $ws.Range("E31") "  +3.55%  "

Set-TextValue $ws.Range("D32") "0.1055"
Set-TextValue $ws.Range("E32") "  +0.86%  "

Set-TextValue $ws.Range("D33") "6.051"
Set-TextValue $ws.Range("E33") "  -2.04%  "

Set-TextValue $ws.Range("D34") "3.763"
Set-TextValue $ws.Range("E34") "  +2.76%  "

Set-TextValue $ws.Range("E35") "  +0.41%  "

Set-TextValue $ws.Range("D36") "0.06549"
Set-TextValue $ws.Range("E36") "  +0.08%  "

Set-TextValue $ws.Range("B37") "Algorand"
Set-TextValue $ws.Range("C37") "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue $ws.Range("D37") "0.2205"
Set-TextValue $ws.Range("E37") "  +0.89%  "

Set-TextValue $ws.Range("B38") "InternetComputer(DFINITY)"
Set-TextValue $ws.Range("C38") "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue $ws.Range("D38") "5.240"
Set-TextValue $ws.Range("E38") "  +2.07%  "

Set-TextValue $ws.Range("D39") "1.218"
Set-TextValue $ws.Range("E39") "  -0.41%  "

Set-TextValue $ws.Range("E40") "  +4.85%  "

Set-TextValue $ws.Range("D41") "0.6500"
Set-TextValue $ws.Range("E41") "  +0.11%  "

Set-TextValue $ws.Range("D42") "8.725"
Set-TextValue $ws.Range("E42") "  -3.35%  "

Set-TextValue $ws.Range("D43") "1.234"
Set-TextValue $ws.Range("E43") "  +0.42%  "

Set-TextValue $ws.Range("D44") "0.6106"
Set-TextValue $ws.Range("E44") "  +1.06%  "

Set-TextValue $ws.Range("D45") "13.25"
Set-TextValue $ws.Range("E45") "  +0.59%  "

Set-TextValue $ws.Range("D46") "3.703"
Set-TextValue $ws.Range("E46") "  +0.76%  "

Set-TextValue $ws.Range("E47") "  +0.59%  "

Set-TextValue $ws.Range("E48") "  +0.58%  "

Set-TextValue $ws.Range("D49") "124.21"
Set-TextValue $ws.Range("E49") "  +0.92%  "

Set-TextValue $ws.Range("D50") "1.165"
Set-TextValue $ws.Range("E50") "  -0.99%  "

Set-TextValue $ws.Range("D51") "78.94"
Set-TextValue $ws.Range("E51") "  +2.13%  "
